$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (item id 5489)
$ws.Range("H2").Value = 339.83334
$ws.Range("I2").Value = 343.54544
$ws.Range("K2").Value = 343.54544
$ws.Range("M2").Value = -230.54544
# Row 9 (item id 5487)
$ws.Range("H9").Value = 5124.9443
$ws.Range("I9").Value = 6126.7334
$ws.Range("J9").Value = 116
$ws.Range("K9").Value = 6126.7334
$ws.Range("L9").Value = 116
$ws.Range("M9").Value = -5957.7334
$ws.Range("N9").Value = -454
# Row 29 (item id 4575)
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 31 (item id 4576)
$ws.Range("H31").Value = 1003
$ws.Range("I31").Value = 1003
$ws.Range("K31").Value = 3009
$ws.Range("M31").Value = -2779
# Row 33 (item id 5512)
$ws.Range("H33").Value = 8333547.5
$ws.Range("I33").Value = 8333547.5
$ws.Range("K33").Value = 8333547.5
$ws.Range("M33").Value = -8333318.5
# Row 41 (item id 5478)
$ws.Range("H41").Value = 1882.6666
$ws.Range("I41").Value = 2575.3333
$ws.Range("J41").Value = 497.33334
$ws.Range("K41").Value = 2575.3333
$ws.Range("L41").Value = 497.33334
$ws.Range("M41").Value = -2135.3333
$ws.Range("N41").Value = -1377.33334
# Row 42 (item id 4600)
$ws.Range("H42").Value = 187.71428
$ws.Range("I42").Value = 54.5
$ws.Range("J42").Value = 365.33334
$ws.Range("K42").Value = 163.5
$ws.Range("L42").Value = 1096.00002
$ws.Range("M42").Value = 66.5
$ws.Range("N42").Value = -1556.00002
# Row 43 (item id 5472)
$ws.Range("H43").Value = 33297.668
$ws.Range("I43").Value = 32409.666
$ws.Range("J43").Value = 35073.668
$ws.Range("K43").Value = 32409.666
$ws.Range("L43").Value = 35073.668
$ws.Range("M43").Value = -32340.666
$ws.Range("N43").Value = -35211.668
# Row 51 (item id 5486)
$ws.Range("H51").Value = 3169.5652
$ws.Range("J51").Value = 4362.5
$ws.Range("L51").Value = 4362.5
$ws.Range("N51").Value = -5330.5
# Row 53 (item id 5479)
$ws.Range("H53").Value = 289.25
$ws.Range("I53").Value = 237.25
$ws.Range("J53").Value = 341.25
$ws.Range("K53").Value = 237.25
$ws.Range("L53").Value = 341.25
$ws.Range("M53").Value = 399.75
$ws.Range("N53").Value = -1615.25
# Row 54 (item id 2174)
$ws.Range("H54").Value = 8522.200000000001
$ws.Range("I54").Value = 7652.75
$ws.Range("J54").Value = 12000
$ws.Range("K54").Value = 7652.75
$ws.Range("L54").Value = 12000
$ws.Range("M54").Value = -7166.75
$ws.Range("N54").Value = -12972
# Row 58 (item id 4606)
$ws.Range("H58").Value = 547.7619
$ws.Range("I58").Value = 305.72223
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 917.16669
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -767.16669
$ws.Range("N58").Value = -6300
# Row 59 (item id 4586)
$ws.Range("H59").Value = 1200
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1200
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 3600
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -4714
# Row 70 (item id 12604)
$ws.Range("H70").Value = 4562.125
$ws.Range("J70").Value = 4562.125
$ws.Range("L70").Value = 13686.375
$ws.Range("N70").Value = -14226.375
# Row 73 (item id 12604)
$ws.Range("H73").Value = 4562.125
$ws.Range("J73").Value = 4562.125
$ws.Range("L73").Value = 13686.375
$ws.Range("N73").Value = -15558.375
# Row 82 (item id 12623)
$ws.Range("H82").Value = 1098
$ws.Range("I82").Value = 1098
$ws.Range("K82").Value = 3294
$ws.Range("M82").Value = -2888
# Row 85 (item id 12623)
$ws.Range("H85").Value = 1098
$ws.Range("I85").Value = 1098
$ws.Range("K85").Value = 3294
$ws.Range("M85").Value = -1890
# Row 87 (item id 10651)
$ws.Range("H87").Value = 81166.5
$ws.Range("J87").Value = 81400
$ws.Range("L87").Value = 81400
$ws.Range("N87").Value = -83896
# Row 90 (item id 10651)
$ws.Range("H90").Value = 81166.5
$ws.Range("J90").Value = 81400
$ws.Range("L90").Value = 244200
$ws.Range("N90").Value = -256680
# Row 97 (item id 19885)
$ws.Range("H97").Value = 222553
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 222553
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 667659
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -668651
# Row 99 (item id 19883)
$ws.Range("H99").Value = 3206.3333
$ws.Range("I99").Value = 447.8
$ws.Range("J99").Value = 6654.5
$ws.Range("K99").Value = 1343.4
$ws.Range("L99").Value = 19963.5
$ws.Range("M99").Value = 154.5999999999999
$ws.Range("N99").Value = -22959.5
# Row 101 (item id 19884)
$ws.Range("H101").Value = 272.66666
$ws.Range("I101").Value = 272.66666
$ws.Range("K101").Value = 817.9999799999999
$ws.Range("M101").Value = 804.0000200000001
# Row 132 (item id 44049)
$ws.Range("H132").Value = 56689
$ws.Range("I132").Value = 31363.94
$ws.Range("K132").Value = 94091.81999999999
$ws.Range("M132").Value = -91561.81999999999
# Row 135 (item id 44047)
$ws.Range("H135").Value = 2081.375
$ws.Range("I135").Value = 1431.6
$ws.Range("J135").Value = 3164.3333
$ws.Range("K135").Value = 12884.4
$ws.Range("L135").Value = 28478.9997
$ws.Range("M135").Value = -10349.4
$ws.Range("N135").Value = -33548.9997
# Row 141 (item id 44161)
$ws.Range("H141").Value = 1750.5
$ws.Range("I141").Value = 1739
$ws.Range("K141").Value = 5217
$ws.Range("M141").Value = -37

$ws = $wb.Worksheets.Item("ARM")
# Row 22 (item id 2497)
$ws.Range("H22").Value = 12964.571
$ws.Range("I22").Value = 858
$ws.Range("J22").Value = 17807.2
$ws.Range("K22").Value = 858
$ws.Range("L22").Value = 17807.2
$ws.Range("M22").Value = -559
$ws.Range("N22").Value = -18405.2
# Row 32 (item id 44147)
$ws.Range("H32").Value = 5716.276
$ws.Range("I32").Value = 5258.88
$ws.Range("K32").Value = 5258.88
$ws.Range("M32").Value = -4971.88
# Row 132 (item id 43997)
$ws.Range("H132").Value = 2493.195
$ws.Range("I132").Value = 2027.2941
$ws.Range("J132").Value = 4756.143
$ws.Range("K132").Value = 6081.8823
$ws.Range("L132").Value = 14268.429
$ws.Range("M132").Value = -3551.8823
$ws.Range("N132").Value = -19328.429

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (item id 5092)
$ws.Range("H22").Value = 224619.89
$ws.Range("I22").Value = 358.27274
$ws.Range("J22").Value = 361668.66
$ws.Range("K22").Value = 358.27274
$ws.Range("L22").Value = 361668.66
$ws.Range("M22").Value = -185.27274
$ws.Range("N22").Value = -362014.66
# Row 132 (item id 41855)
$ws.Range("H132").Value = 82489.5
$ws.Range("J132").Value = 82489.5
$ws.Range("L132").Value = 82489.5
$ws.Range("N132").Value = -92609.5
# Row 134 (item id 43998)
$ws.Range("H134").Value = 2001.9387
$ws.Range("I134").Value = 1784.674
$ws.Range("J134").Value = 5333.3335
$ws.Range("K134").Value = 5354.022
$ws.Range("L134").Value = 16000.0005
$ws.Range("M134").Value = -2819.022
$ws.Range("N134").Value = -21070.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 111 (item id 25792)
$ws.Range("H111").Value = 45000
$ws.Range("J111").Value = 45000
$ws.Range("L111").Value = 45000
$ws.Range("N111").Value = -53180
# Row 122 (item id 36196)
$ws.Range("H122").Value = 1373.3846
$ws.Range("I122").Value = 1236.6666
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 3709.9998
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -1259.9998
$ws.Range("N122").Value = -13942
# Row 134 (item id 44020)
$ws.Range("H134").Value = 14651.514
$ws.Range("I134").Value = 8428.700000000001
$ws.Range("K134").Value = 25286.1
$ws.Range("M134").Value = -22751.1

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (item id 36169)
$ws.Range("H102").Value = 1795.238
$ws.Range("I102").Value = 1821.6111
$ws.Range("J102").Value = 1637
$ws.Range("K102").Value = 1821.6111
$ws.Range("L102").Value = 1637
$ws.Range("M102").Value = -199.6111000000001
$ws.Range("N102").Value = -4881

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (item id 5282)
$ws.Range("H46").Value = 1354.3077
$ws.Range("I46").Value = 1371.8
$ws.Range("J46").Value = 1296
$ws.Range("K46").Value = 1371.8
$ws.Range("L46").Value = 1296
$ws.Range("M46").Value = -1183.8
$ws.Range("N46").Value = -1672
# Row 55 (item id 5284)
$ws.Range("H55").Value = 322.78946
$ws.Range("I55").Value = 385.1111
$ws.Range("J55").Value = 266.7
$ws.Range("K55").Value = 385.1111
$ws.Range("L55").Value = 266.7
$ws.Range("M55").Value = -212.1111
$ws.Range("N55").Value = -612.7
# Row 132 (item id 44058)
$ws.Range("H132").Value = 5179
$ws.Range("I132").Value = 4245.5
$ws.Range("K132").Value = 12736.5
$ws.Range("M132").Value = -10206.5
# Row 136 (item id 44060)
$ws.Range("H136").Value = 7212.4287
$ws.Range("I136").Value = 6663.778
$ws.Range("K136").Value = 19991.334
$ws.Range("M136").Value = -17441.334

$ws = $wb.Worksheets.Item("WVR")
# Row 4 (item id 2996)
$ws.Range("H4").Value = 3887.0833
$ws.Range("I4").Value = 3737.25
$ws.Range("J4").Value = 3962
$ws.Range("K4").Value = 3737.25
$ws.Range("L4").Value = 3962
$ws.Range("M4").Value = -3624.25
$ws.Range("N4").Value = -4188
# Row 62 (item id 12589)
$ws.Range("H62").Value = 17358.938
$ws.Range("J62").Value = 9033.166999999999
$ws.Range("L62").Value = 9033.166999999999
$ws.Range("N62").Value = -10281.167
# Row 65 (item id 12589)
$ws.Range("H65").Value = 17358.938
$ws.Range("J65").Value = 9033.166999999999
$ws.Range("L65").Value = 45165.835
$ws.Range("N65").Value = -51405.835
# Row 75 (item id 11957)
$ws.Range("H75").Value = 39875
$ws.Range("J75").Value = 39875
$ws.Range("L75").Value = 39875
$ws.Range("N75").Value = -41747
# Row 78 (item id 11957)
$ws.Range("H78").Value = 39875
$ws.Range("J78").Value = 39875
$ws.Range("L78").Value = 119625
$ws.Range("N78").Value = -128985
# Row 105 (item id 18710)
$ws.Range("H105").Value = 31665
$ws.Range("J105").Value = 31665
$ws.Range("L105").Value = 31665
$ws.Range("N105").Value = -38653
# Row 113 (item id 27752)
$ws.Range("H113").Value = 762.25
$ws.Range("J113").Value = 766.6667
$ws.Range("L113").Value = 2300.0001
$ws.Range("N113").Value = -6640.0001

$wb.Save()